$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 06:44:48"
$wsZhCn.Range("H2").Value = "2016-03-12 06:45:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 06:44:51"
$wsDeDe.Range("H2").Value = "2016-03-12 06:45:10"
